$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.511.05"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "3.361.46"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "560.97"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "176.20"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("D8").Value = "3.342.71"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.633"
$ws.Range("E10").Value = "  +3.96%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").Value = "  +10.30%  "
$ws.Range("D12").Value = "55.80"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").Value = "9.12"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "3.889.67"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "18.32"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.118"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.356.92"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "11.86"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "64.362.83"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "0.992"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "463.25"
$ws.Range("E22").Value = "  +15.67%  "
$ws.Range("D23").Value = "4.89"
$ws.Range("E23").Value = "  +11.26%  "
$ws.Range("D24").Value = "4.10"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "86.33"
$ws.Range("E25").Value = "  +5.23%  "
$ws.Range("D26").Value = "13.58"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("D29").Value = "8.83"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "30.23"
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("D31").Value = "6.71"
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("D32").Value = "11.50"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "579.73"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "59.17"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -5.45%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "36.00"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "3.092.38"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "8.37"
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "136.48"
$ws.Range("E51").Value = "  +1.69%  "
